# Generate Report for Handoff
# Adds a new tracked file (9e9f81f3-54e6-47e2-9ea6-80c598a4ef0a.md) as row 3
# on the "Overview", "zh-cn" and "de-de" sheets, mirroring the existing
# 57f1b5db-... row.
#
# NOTE: a leading "'" forces Excel to store the value as literal text
# (matching the source workbook, which stores "True"/"False"/"" as text,
# not as booleans / empty cells) without the apostrophe itself being kept.

$wb = $excel.ActiveWorkbook

$guid = "9e9f81f3-54e6-47e2-9ea6-80c598a4ef0a"
$blobUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/34fb509076e49b457389259e65d9400024201c2a/e2e/" + $guid + ".md"

# ---------------------------------------------------------------------------
# Overview sheet (row 3)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A3").Value = $guid + ".md"

$ws1.Range("B3").Value = "e2e\" + $guid + ".md"
$ws1.Hyperlinks.Add($ws1.Range("B3"), $blobUrl, "", "", "e2e\" + $guid + ".md") | Out-Null
$ws1.Range("B3").Style = "HyperLink"

$ws1.Range("C3").Value = ".md"
$ws1.Range("D3").Value = "'"
$ws1.Range("E3").Value = "Ready for handoff"
$ws1.Range("F3").Value = "Ready for handoff"
$ws1.Range("G3").Value = "2016-08-28 14:39:32"
$ws1.Range("G3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$ws1.ListObjects.Item(1).Resize($ws1.Range("A1:G3"))

# ---------------------------------------------------------------------------
# zh-cn sheet (row 3)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A3").Value = $guid + ".md"
$ws2.Hyperlinks.Add($ws2.Range("A3"), $blobUrl, "", "", $guid + ".md") | Out-Null
$ws2.Range("A3").Style = "HyperLink"

$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = "Ready for handoff"
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "'False"
$ws2.Range("G3").Value = $guid + ".784865fd1f6d0f5a44cbd9e64332fe42c063a2c6.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-28 14:39:27"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("I3").Value = "'"
$ws2.Range("J3").Value = "'"
$ws2.Range("K3").Value = "0001-01-01 00:00:00"
$ws2.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("L3").Value = "'"
$ws2.Range("M3").Value = "'True"
$ws2.Range("N3").Value = "'"
$ws2.Range("O3").Value = "'False"
$ws2.Range("P3").Value = "'"

$ws2.ListObjects.Item(1).Resize($ws2.Range("A1:P3"))

# ---------------------------------------------------------------------------
# de-de sheet (row 3)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A3").Value = $guid + ".md"
$ws3.Hyperlinks.Add($ws3.Range("A3"), $blobUrl, "", "", $guid + ".md") | Out-Null
$ws3.Range("A3").Style = "HyperLink"

$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = "Ready for handoff"
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "'False"
$ws3.Range("G3").Value = $guid + ".784865fd1f6d0f5a44cbd9e64332fe42c063a2c6.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-28 14:39:32"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("I3").Value = "'"
$ws3.Range("J3").Value = "'"
$ws3.Range("K3").Value = "0001-01-01 00:00:00"
$ws3.Range("K3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("L3").Value = "'"
$ws3.Range("M3").Value = "'True"
$ws3.Range("N3").Value = "'"
$ws3.Range("O3").Value = "'False"
$ws3.Range("P3").Value = "'"

$ws3.ListObjects.Item(1).Resize($ws3.Range("A1:P3"))
